$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowFV {
    param($row, $vals)
    $col = 6
    foreach ($v in $vals) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

Set-RowFV 63 @('Villarreal', 1, 'Girona', 2, 1.71, '17/09/2023 09:02', 2.17, '27/09/2023 18:51', 4.14, '17/09/2023 09:02', 3.72, '27/09/2023 18:51', 4.72, '17/09/2023 09:02', 3.42, '27/09/2023 18:51', 'https://www.betexplorer.com/football/spain/laliga/villarreal-girona/80EuTg3A/')
Set-RowFV 64 @('Real Madrid', 2, 'Las Palmas', 0, 1.2, '23/09/2023 09:28', 1.18, '27/09/2023 18:29', 6.76, '23/09/2023 09:28', 8, '27/09/2023 18:29', 11.3, '23/09/2023 09:28', 16.5, '27/09/2023 18:29', 'https://www.betexplorer.com/football/spain/laliga/real-madrid-las-palmas/GQHmRXXM/')
Set-RowFV 65 @('Ath Bilbao', 2, 'Getafe', 2, 1.71, '17/09/2023 09:02', 1.53, '27/09/2023 18:31', 3.42, '17/09/2023 09:02', 4.06, '27/09/2023 18:49', 5.44, '17/09/2023 09:02', 7.73, '27/09/2023 18:49', 'https://www.betexplorer.com/football/spain/laliga/ath-bilbao-getafe/zgsFCYIT/')
Set-RowFV 66 @('Valencia', 0, 'Real Sociedad', 1, 2.33, '17/09/2023 09:02', 2.57, '27/09/2023 21:27', 3.14, '17/09/2023 09:02', 3.02, '27/09/2023 21:27', 3.26, '17/09/2023 09:02', 3.31, '27/09/2023 21:27', 'https://www.betexplorer.com/football/spain/laliga/valencia-real-sociedad/M3IqSDIG/')
Set-RowFV 67 @('Cadiz CF', 0, 'Rayo Vallecano', 0, 2.59, '17/09/2023 09:02', 2.72, '27/09/2023 21:19', 3.06, '17/09/2023 09:02', 3.11, '27/09/2023 21:17', 3.11, '17/09/2023 09:02', 3, '27/09/2023 21:30', 'https://www.betexplorer.com/football/spain/laliga/cadiz-rayo-vallecano/CEYt8hRp/')
Set-RowFV 88 @('Alaves', 1, 'Betis', 1, 2.8, '01/10/2023 20:24', 2.57, '08/10/2023 18:28', 3.04, '01/10/2023 20:24', 3.25, '08/10/2023 18:28', 2.72, '01/10/2023 20:24', 3.06, '08/10/2023 18:22', 'https://www.betexplorer.com/football/spain/laliga/alaves-betis/YNPlfW19/')
Set-RowFV 89 @('Celta Vigo', 2, 'Getafe', 2, 1.92, '28/09/2023 15:02', 2.04, '08/10/2023 18:29', 3.25, '28/09/2023 15:02', 3.29, '08/10/2023 18:27', 4.85, '28/09/2023 15:02', 4.38, '08/10/2023 18:29', 'https://www.betexplorer.com/football/spain/laliga/celta-vigo-getafe/0ARtdhXd/')
Set-RowFV 172 @('Atl. Madrid', 3, 'Getafe', 3, 1.31, '10/12/2023 00:01', 1.48, '19/12/2023 21:28', 4.91, '10/12/2023 00:01', 4.54, '19/12/2023 21:29', 9, '10/12/2023 00:01', 7.35, '19/12/2023 21:29', 'https://www.betexplorer.com/football/spain/laliga/atl-madrid-getafe/EwmSwQZI/')
Set-RowFV 173 @('Granada CF', 0, 'Sevilla', 3, 2.76, '10/12/2023 00:01', 3.12, '19/12/2023 21:27', 3.33, '10/12/2023 00:01', 3.17, '19/12/2023 21:05', 2.48, '10/12/2023 00:01', 2.55, '19/12/2023 21:25', 'https://www.betexplorer.com/football/spain/laliga/granada-cf-sevilla/KhVmeR4g/')

function Set-RowAV {
    param($row, $vals)
    $col = 1
    foreach ($v in $vals) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

# Append two new rows (177, 178), copying formatting from row 176 first
$ws.Range("A176:V176").Copy() | Out-Null
$ws.Range("A177:V178").PasteSpecial(-4122) | Out-Null

Set-RowAV 177 @(176, 'spain', 'laliga', '2023-2024', 45281.79166666666, 'Betis', 1, 'Girona', 1, 2.31, '10/12/2023 10:02', 2.67, '21/12/2023 18:59', 3.42, '10/12/2023 10:02', 3.55, '21/12/2023 18:59', 2.96, '10/12/2023 10:02', 2.69, '21/12/2023 18:58', 'https://www.betexplorer.com/football/spain/laliga/betis-girona/Sfj929CJ/')
Set-RowAV 178 @(177, 'spain', 'laliga', '2023-2024', 45281.79166666666, 'Cadiz CF', 0, 'Real Sociedad', 0, 4.08, '10/12/2023 10:02', 6.27, '21/12/2023 18:57', 3.38, '10/12/2023 10:02', 3.59, '21/12/2023 18:59', 1.9, '10/12/2023 10:02', 1.69, '21/12/2023 18:58', 'https://www.betexplorer.com/football/spain/laliga/cadiz-real-sociedad/6okD1TRP/')